$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.685.92'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.73'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  -2.86%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.18'
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("E6").Value = '  -2.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4311'
$ws.Range("E7").Value = '  -2.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3752'
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07352'
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8810'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.858.11'
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.737'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.459'
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07104'
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.73'
$ws.Range("E16").Value = '  +4.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.014'
$ws.Range("E17").Value = '  -2.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008994'
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.011'
$ws.Range("E19").Value = '  -2.64%  '
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.694.76'
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.258'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.18'
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.079.05'
$ws.Range("E24").Value = '  -1.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.037'
$ws.Range("E25").Value = '  +2.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.57'
$ws.Range("E26").Value = '  -2.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.60'
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.141'
$ws.Range("E28").Value = '  +7.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.389'
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.16'
$ws.Range("E30").Value = '  +1.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08921'
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.230'
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7775'
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.560'
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.911'
$ws.Range("E35").Value = '  -6.63%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.012'
$ws.Range("E36").Value = '  -2.65%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.140'
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05343'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01972'
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.220'
$ws.Range("E40").Value = '  +4.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.863'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5172'
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1679'
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.955'
$ws.Range("E44").Value = '  +3.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.72'
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4735'
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.703'
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06499'
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.012'
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.891'
$ws.Range("E51").Value = '  -0.90%  '
